$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "keshav"
$ws.Range("B3").Value = "keshav"

$ws.Range("B3").Select()
